$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 3 (which holds 252980),
# shifting it down to row 5.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Update the value in row 2 (was 252417, now 252466)
$ws.Range("A2").Value = 252466

# Fill in the two newly inserted rows
$ws.Range("A3").Value = 252417
$ws.Range("B3").Value = "nessuna compatibilità con alcuna macchina"

$ws.Range("A4").Value = 252418
$ws.Range("B4").Value = "nessuna compatibilità con alcuna macchina"
